$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swerve Home table (columns E:G) ---

# Header row
$ws.Range("E1").Value = "Swerve Home"
$ws.Range("F1").Value = "Home Relative to Halsensor"
$ws.Range("G1").Value = "Halsensor position relative to home"

# Front left
$ws.Range("E2").Value = "Front left"
$ws.Range("F2").Value = -8.9524135589599592
$ws.Range("G2").Formula = "=-F2+9"

# Front right (shared formula master)
$ws.Range("E3").Value = "Front right"
$ws.Range("F3").Value = -8.7857446670532209
$ws.Range("G3:G5").Formula = "=-F3+9"

# Back left
$ws.Range("E4").Value = "Back left"
$ws.Range("F4").Value = -8.9047937393188406

# Back right
$ws.Range("E5").Value = "Back right"
$ws.Range("F5").Value = -8.9524135589599592

# --- Swerve drive PID table (columns E:F) ---

$ws.Range("E8").Value = "Swerve drive PID"

$ws.Range("E9").Value = "P"
$ws.Range("F9").Value = 0.00001

$ws.Range("E10").Value = "I"

$ws.Range("E11").Value = "D"

$ws.Range("E12").Value = "F"
$ws.Range("F12").Value = 0.000166

# --- Column widths for the new columns ---
$ws.Columns.Item(5).ColumnWidth = 11.4
$ws.Columns.Item(6).ColumnWidth = 23.1
$ws.Columns.Item(7).ColumnWidth = 30

# --- Selection as left by the author ---
$ws.Range("H14").Select() | Out-Null
